$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.1209964412811388
$wsSummary.Range("C2").Value = 0.05363984674329502
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.1018181818181818
$wsSummary.Range("F2").Value = 0.2208201892744479
$wsSummary.Range("G2").Value = 0.5957446808510638
$wsSummary.Range("H2").Value = 0.7077314071696094
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 494
$wsSummary.Range("K2").Value = 40
$wsSummary.Range("L2").Value = 0

# --- Sheet 2: Classification Report ---
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("B2").Value = 1
$wsReport.Range("C2").Value = 0.0749063670411985
$wsReport.Range("D2").Value = 0.1393728222996516

$wsReport.Range("B3").Value = 0.05363984674329502
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.1018181818181818

$wsReport.Range("B4").Value = 0.1209964412811388
$wsReport.Range("C4").Value = 0.1209964412811388
$wsReport.Range("D4").Value = 0.1209964412811388
$wsReport.Range("E4").Value = 0.1209964412811388

$wsReport.Range("B5").Value = 0.5268199233716475
$wsReport.Range("C5").Value = 0.5374531835205992
$wsReport.Range("D5").Value = 0.1205955020589167

$wsReport.Range("B6").Value = 0.9528503838235093
$wsReport.Range("C6").Value = 0.1209964412811388
$wsReport.Range("D6").Value = 0.1375017725959485

# --- Sheet 3: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 40
$wsConf.Range("C2").Value = 494
$wsConf.Range("B3").Value = 0
$wsConf.Range("C3").Value = 28
